$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure columns -------------------------------------------------
# Original layout: A Group_Name, B Group_Type_ID, C Ministry_ID,
#   D Congregation_ID, E Primary_Contact, F Start_Date, G Target_Size,
#   H Description, (J sparse free-text notes)
#
# Target layout: A Group_Name, B Group_Type_ID, C Ministry_ID,
#   D Congregation_ID, E Primary_Contact, F Is_Primary_Contact_Host,
#   G Start_Date, H Enable_Waiting_List, I Child_Care_Available,
#   J Meeting_Time, K Target_Size, L Description,
#   M Deadline_Passed_Message_ID, N __IsPublic, O __ISBlogEnabled,
#   P __ISWebEnabled

# Remove the old free-text notes column first. (Clearing contents rather
# than deleting the column avoids shifting/duplicating its custom-width
# metadata onto a neighbouring column.)
$ws.Range("J1:J19").ClearContents()

# Insert the new "Is_Primary_Contact_Host" column before the old Start_Date.
$ws.Columns("F:F").Insert()

# Insert three new columns before the old Target_Size column (now column H).
$ws.Columns("H:J").Insert()

# The inserted columns pick up the neighbouring date-format styling; drop it
# so untouched rows don't leave behind empty-but-styled cells.
$ws.Range("H2:J19").Clear()

# --- Header row (order chosen to mirror the original authoring sequence,
#     which drives the order new entries land in the shared-strings table)
$ws.Range("M1").Value2 = "Deadline_Passed_Message_ID"
$ws.Range("N1").Value2 = "__IsPublic"
$ws.Range("O1").Value2 = "__ISBlogEnabled"
$ws.Range("P1").Value2 = "__ISWebEnabled"

$ws.Range("H1").Value2 = "Enable_Waiting_List"
$ws.Range("I1").Value2 = "Child_Care_Available"
$ws.Range("J1").Value2 = "Meeting_Time"

# --- Misc new per-row values -----------------------------------------------
$ws.Range("H6").Value2 = 1          # Enable_Waiting_List
$ws.Range("I15").Value2 = 1         # Child_Care_Available

$ws.Range("J17").Value2 = 0.70833333333333337   # Meeting_Time
$ws.Range("J17").NumberFormat = "h:mm:ss"

$deadlineRows = @(6,7,8,9,13,14,16,17,18,19)
foreach ($r in $deadlineRows) {
  $ws.Cells.Item($r, 13).Value2 = 58
}
$ws.Range("M15").Value2 = 59

$publicRows = @(2,3,8,9,12)
foreach ($r in $publicRows) {
  $ws.Cells.Item($r, 14).Value2 = "N"
  $ws.Cells.Item($r, 15).Value2 = "Y"
  $ws.Cells.Item($r, 16).Value2 = "Y"
}

# Header typed last of all the new text labels.
$ws.Range("F1").Value2 = "Is_Primary_Contact_Host"

# --- Is_Primary_Contact_Host values (column F) for every data row --------
$hostValues = @{
  2 = 0; 3 = 0; 4 = 0; 5 = 0; 6 = 0; 7 = 0; 8 = 0; 9 = 0; 10 = 0; 11 = 0;
  12 = 0; 13 = 0; 14 = 0; 15 = 0; 16 = 1; 17 = 1; 18 = 1; 19 = 1
}
foreach ($r in $hostValues.Keys) {
  $ws.Cells.Item($r, 6).Value2 = $hostValues[$r]
}

# --- Column widths (best effort match of authored widths) -----------------
$ws.Range("E1").ColumnWidth = 34.8
$ws.Range("F1").ColumnWidth = 20.67
$ws.Range("H1").ColumnWidth = 16.86
$ws.Range("I1").ColumnWidth = 17.4
$ws.Range("J1").ColumnWidth = 12.13
$ws.Range("M1").ColumnWidth = 24.53
$ws.Range("N1").ColumnWidth = 8.53
$ws.Range("O1").ColumnWidth = 13.73
$ws.Range("P1").ColumnWidth = 13.93

# --- Selection matches the authored state ---------------------------------
$ws.Range("H25").Select()
